$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17: Directory column (D) used lowercase "data/xml/" (was "data/XML/")
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("D$r").Value = "data/xml/"
}

# Rows 18-35: Directory column (D) used lowercase "nodegoat/" (was "Nodegoat/")
for ($r = 18; $r -le 35; $r++) {
    $ws.Range("D$r").Value = "nodegoat/"
}

# New explicit width for column D
$ws.Columns("D:D").ColumnWidth = 17.4

# Update the active selection
$ws.Range("E29").Select()
